$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("final_fail")
$ws.Cells.Item(2, 1).Value = "Submissions (% of course total)"
$ws.Cells.Item(2, 2).Value = $true
$ws.Cells.Item(2, 3).Value = $true
$ws.Cells.Item(2, 4).Value = $true
$ws.Cells.Item(2, 5).Value = $true
$ws.Cells.Item(2, 6).Value = $true
$ws.Cells.Item(2, 7).Value = $false
$ws.Cells.Item(2, 8).Value = $true
$ws.Cells.Item(2, 9).Value = $false
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(3, 1).Value = "Total time online (min)"
$ws.Cells.Item(3, 2).Value = $true
$ws.Cells.Item(3, 3).Value = $true
$ws.Cells.Item(3, 4).Value = $true
$ws.Cells.Item(3, 5).Value = $true
$ws.Cells.Item(3, 6).Value = $true
$ws.Cells.Item(3, 7).Value = $false
$ws.Cells.Item(3, 8).Value = $true
$ws.Cells.Item(3, 9).Value = $false
$ws.Cells.Item(3, 10).Value = 6
$ws.Cells.Item(4, 1).Value = "On/off campus click ratio"
$ws.Cells.Item(4, 2).Value = $true
$ws.Cells.Item(4, 3).Value = $true
$ws.Cells.Item(4, 4).Value = $true
$ws.Cells.Item(4, 5).Value = $true
$ws.Cells.Item(4, 6).Value = $true
$ws.Cells.Item(4, 7).Value = $false
$ws.Cells.Item(4, 8).Value = $true
$ws.Cells.Item(4, 9).Value = $false
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(5, 1).Value = "Clicks (% of course total)"
$ws.Cells.Item(5, 2).Value = $true
$ws.Cells.Item(5, 3).Value = $true
$ws.Cells.Item(5, 4).Value = $true
$ws.Cells.Item(5, 5).Value = $true
$ws.Cells.Item(5, 6).Value = $true
$ws.Cells.Item(5, 7).Value = $false
$ws.Cells.Item(5, 8).Value = $true
$ws.Cells.Item(5, 9).Value = $false
$ws.Cells.Item(5, 10).Value = 6
$ws.Cells.Item(6, 1).Value = "Resources viewed"
$ws.Cells.Item(6, 2).Value = $true
$ws.Cells.Item(6, 3).Value = $true
$ws.Cells.Item(6, 4).Value = $true
$ws.Cells.Item(6, 5).Value = $true
$ws.Cells.Item(6, 6).Value = $true
$ws.Cells.Item(6, 7).Value = $false
$ws.Cells.Item(6, 8).Value = $true
$ws.Cells.Item(6, 9).Value = $false
$ws.Cells.Item(6, 10).Value = 6
$ws.Cells.Item(7, 1).Value = "Largest period of inactivity (h)"
$ws.Cells.Item(7, 2).Value = $true
$ws.Cells.Item(7, 3).Value = $true
$ws.Cells.Item(7, 4).Value = $false
$ws.Cells.Item(7, 5).Value = $true
$ws.Cells.Item(7, 6).Value = $true
$ws.Cells.Item(7, 7).Value = $true
$ws.Cells.Item(7, 8).Value = $true
$ws.Cells.Item(7, 9).Value = $false
$ws.Cells.Item(7, 10).Value = 6
$ws.Cells.Item(8, 1).Value = "Clicks on campus"
$ws.Cells.Item(8, 2).Value = $true
$ws.Cells.Item(8, 3).Value = $true
$ws.Cells.Item(8, 4).Value = $true
$ws.Cells.Item(8, 5).Value = $true
$ws.Cells.Item(8, 6).Value = $true
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = $true
$ws.Cells.Item(8, 9).Value = $false
$ws.Cells.Item(8, 10).Value = 6
$ws.Cells.Item(9, 1).Value = "Average session duration (min)"
$ws.Cells.Item(9, 2).Value = $true
$ws.Cells.Item(9, 3).Value = $true
$ws.Cells.Item(9, 4).Value = $false
$ws.Cells.Item(9, 5).Value = $true
$ws.Cells.Item(9, 6).Value = $true
$ws.Cells.Item(9, 7).Value = $false
$ws.Cells.Item(9, 8).Value = $true
$ws.Cells.Item(9, 9).Value = $false
$ws.Cells.Item(9, 10).Value = 5
$ws.Cells.Item(10, 1).Value = "Number of days"
$ws.Cells.Item(10, 2).Value = $true
$ws.Cells.Item(10, 3).Value = $true
$ws.Cells.Item(10, 4).Value = $false
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(10, 6).Value = $true
$ws.Cells.Item(10, 7).Value = $false
$ws.Cells.Item(10, 8).Value = $true
$ws.Cells.Item(10, 9).Value = $false
$ws.Cells.Item(10, 10).Value = 5
$ws.Cells.Item(11, 1).Value = "Clicks per day"
$ws.Cells.Item(11, 2).Value = $true
$ws.Cells.Item(11, 3).Value = $true
$ws.Cells.Item(11, 4).Value = $false
$ws.Cells.Item(11, 5).Value = $true
$ws.Cells.Item(11, 6).Value = $true
$ws.Cells.Item(11, 7).Value = $false
$ws.Cells.Item(11, 8).Value = $true
$ws.Cells.Item(11, 9).Value = $false
$ws.Cells.Item(11, 10).Value = 5
$ws.Cells.Item(12, 1).Value = "Start of Session 1 (%)"
$ws.Cells.Item(12, 2).Value = $true
$ws.Cells.Item(12, 3).Value = $true
$ws.Cells.Item(12, 4).Value = $false
$ws.Cells.Item(12, 5).Value = $true
$ws.Cells.Item(12, 6).Value = $true
$ws.Cells.Item(12, 7).Value = $false
$ws.Cells.Item(12, 8).Value = $true
$ws.Cells.Item(12, 9).Value = $false
$ws.Cells.Item(12, 10).Value = 5
$ws.Cells.Item(13, 1).Value = "Clicks per session"
$ws.Cells.Item(13, 2).Value = $true
$ws.Cells.Item(13, 3).Value = $true
$ws.Cells.Item(13, 4).Value = $false
$ws.Cells.Item(13, 5).Value = $true
$ws.Cells.Item(13, 6).Value = $true
$ws.Cells.Item(13, 7).Value = $false
$ws.Cells.Item(13, 8).Value = $true
$ws.Cells.Item(13, 9).Value = $false
$ws.Cells.Item(13, 10).Value = 5
$ws.Cells.Item(14, 1).Value = "Days with no interaction (%)"
$ws.Cells.Item(14, 2).Value = $true
$ws.Cells.Item(14, 3).Value = $true
$ws.Cells.Item(14, 4).Value = $false
$ws.Cells.Item(14, 5).Value = $true
$ws.Cells.Item(14, 6).Value = $true
$ws.Cells.Item(14, 7).Value = $false
$ws.Cells.Item(14, 8).Value = $true
$ws.Cells.Item(14, 9).Value = $false
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(15, 1).Value = "Average grade of assignments"
$ws.Cells.Item(15, 2).Value = $true
$ws.Cells.Item(15, 3).Value = $true
$ws.Cells.Item(15, 4).Value = $true
$ws.Cells.Item(15, 5).Value = $true
$ws.Cells.Item(15, 6).Value = $false
$ws.Cells.Item(15, 7).Value = $false
$ws.Cells.Item(15, 8).Value = $true
$ws.Cells.Item(15, 9).Value = $false
$ws.Cells.Item(15, 10).Value = 5
$ws.Cells.Item(16, 1).Value = "Start of Session 4 (%)"
$ws.Cells.Item(16, 2).Value = $true
$ws.Cells.Item(16, 3).Value = $true
$ws.Cells.Item(16, 4).Value = $false
$ws.Cells.Item(16, 5).Value = $false
$ws.Cells.Item(16, 6).Value = $true
$ws.Cells.Item(16, 7).Value = $false
$ws.Cells.Item(16, 8).Value = $true
$ws.Cells.Item(16, 9).Value = $false
$ws.Cells.Item(16, 10).Value = 4
$ws.Cells.Item(17, 1).Value = "Days with no interaction"
$ws.Cells.Item(17, 2).Value = $true
$ws.Cells.Item(17, 3).Value = $false
$ws.Cells.Item(17, 4).Value = $false
$ws.Cells.Item(17, 5).Value = $true
$ws.Cells.Item(17, 6).Value = $true
$ws.Cells.Item(17, 7).Value = $false
$ws.Cells.Item(17, 8).Value = $true
$ws.Cells.Item(17, 9).Value = $false
$ws.Cells.Item(17, 10).Value = 4
$ws.Cells.Item(18, 1).Value = "Start of Session 2 (%)"
$ws.Cells.Item(18, 2).Value = $false
$ws.Cells.Item(18, 3).Value = $true
$ws.Cells.Item(18, 4).Value = $false
$ws.Cells.Item(18, 5).Value = $true
$ws.Cells.Item(18, 6).Value = $true
$ws.Cells.Item(18, 7).Value = $false
$ws.Cells.Item(18, 8).Value = $true
$ws.Cells.Item(18, 9).Value = $false
$ws.Cells.Item(18, 10).Value = 4
$ws.Cells.Item(19, 1).Value = "Clicks on course"
$ws.Cells.Item(19, 2).Value = $true
$ws.Cells.Item(19, 3).Value = $false
$ws.Cells.Item(19, 4).Value = $false
$ws.Cells.Item(19, 5).Value = $false
$ws.Cells.Item(19, 6).Value = $true
$ws.Cells.Item(19, 7).Value = $false
$ws.Cells.Item(19, 8).Value = $true
$ws.Cells.Item(19, 9).Value = $false
$ws.Cells.Item(19, 10).Value = 3
$ws.Cells.Item(20, 1).Value = "Assignments viewed"
$ws.Cells.Item(20, 2).Value = $true
$ws.Cells.Item(20, 3).Value = $false
$ws.Cells.Item(20, 4).Value = $false
$ws.Cells.Item(20, 5).Value = $false
$ws.Cells.Item(20, 6).Value = $false
$ws.Cells.Item(20, 7).Value = $false
$ws.Cells.Item(20, 8).Value = $true
$ws.Cells.Item(20, 9).Value = $false
$ws.Cells.Item(20, 10).Value = 2
$ws.Cells.Item(21, 1).Value = "Files downloaded"
$ws.Cells.Item(21, 2).Value = $false
$ws.Cells.Item(21, 3).Value = $false
$ws.Cells.Item(21, 4).Value = $true
$ws.Cells.Item(21, 5).Value = $false
$ws.Cells.Item(21, 6).Value = $false
$ws.Cells.Item(21, 7).Value = $false
$ws.Cells.Item(21, 8).Value = $true
$ws.Cells.Item(21, 9).Value = $false
$ws.Cells.Item(21, 10).Value = 2
$ws.Cells.Item(22, 1).Value = "Forum posts"
$ws.Cells.Item(22, 2).Value = $false
$ws.Cells.Item(22, 3).Value = $false
$ws.Cells.Item(22, 4).Value = $true
$ws.Cells.Item(22, 5).Value = $false
$ws.Cells.Item(22, 6).Value = $false
$ws.Cells.Item(22, 7).Value = $false
$ws.Cells.Item(22, 8).Value = $true
$ws.Cells.Item(22, 9).Value = $false
$ws.Cells.Item(22, 10).Value = 2
$ws.Cells.Item(23, 1).Value = "Number of clicks"
$ws.Cells.Item(23, 2).Value = $false
$ws.Cells.Item(23, 3).Value = $false
$ws.Cells.Item(23, 4).Value = $false
$ws.Cells.Item(23, 5).Value = $false
$ws.Cells.Item(23, 6).Value = $true
$ws.Cells.Item(23, 7).Value = $false
$ws.Cells.Item(23, 8).Value = $true
$ws.Cells.Item(23, 9).Value = $false
$ws.Cells.Item(23, 10).Value = 2
$ws.Cells.Item(24, 1).Value = "Clicks on folder"
$ws.Cells.Item(24, 2).Value = $false
$ws.Cells.Item(24, 3).Value = $false
$ws.Cells.Item(24, 4).Value = $true
$ws.Cells.Item(24, 5).Value = $false
$ws.Cells.Item(24, 6).Value = $false
$ws.Cells.Item(24, 7).Value = $false
$ws.Cells.Item(24, 8).Value = $true
$ws.Cells.Item(24, 9).Value = $false
$ws.Cells.Item(24, 10).Value = 2
$ws.Cells.Item(25, 1).Value = "Assignments submitted"
$ws.Cells.Item(25, 2).Value = $false
$ws.Cells.Item(25, 3).Value = $false
$ws.Cells.Item(25, 4).Value = $true
$ws.Cells.Item(25, 5).Value = $false
$ws.Cells.Item(25, 6).Value = $false
$ws.Cells.Item(25, 7).Value = $false
$ws.Cells.Item(25, 8).Value = $true
$ws.Cells.Item(25, 9).Value = $false
$ws.Cells.Item(25, 10).Value = 2
$ws.Cells.Item(26, 1).Value = "Number of sessions"
$ws.Cells.Item(26, 2).Value = $false
$ws.Cells.Item(26, 3).Value = $false
$ws.Cells.Item(26, 4).Value = $true
$ws.Cells.Item(26, 5).Value = $false
$ws.Cells.Item(26, 6).Value = $false
$ws.Cells.Item(26, 7).Value = $false
$ws.Cells.Item(26, 8).Value = $true
$ws.Cells.Item(26, 9).Value = $false
$ws.Cells.Item(26, 10).Value = 2
$ws.Cells.Item(27, 1).Value = "Start of Session 7 (%)"
$ws.Cells.Item(27, 2).Value = $false
$ws.Cells.Item(27, 3).Value = $false
$ws.Cells.Item(27, 4).Value = $false
$ws.Cells.Item(27, 5).Value = $false
$ws.Cells.Item(27, 6).Value = $true
$ws.Cells.Item(27, 7).Value = $false
$ws.Cells.Item(27, 8).Value = $true
$ws.Cells.Item(27, 9).Value = $false
$ws.Cells.Item(27, 10).Value = 2
$ws.Cells.Item(28, 1).Value = "Start of Session 5 (%)"
$ws.Cells.Item(28, 2).Value = $false
$ws.Cells.Item(28, 3).Value = $false
$ws.Cells.Item(28, 4).Value = $false
$ws.Cells.Item(28, 5).Value = $false
$ws.Cells.Item(28, 6).Value = $true
$ws.Cells.Item(28, 7).Value = $false
$ws.Cells.Item(28, 8).Value = $true
$ws.Cells.Item(28, 9).Value = $false
$ws.Cells.Item(28, 10).Value = 2
$ws.Cells.Item(29, 1).Value = "Start of Session 3 (%)"
$ws.Cells.Item(29, 2).Value = $false
$ws.Cells.Item(29, 3).Value = $false
$ws.Cells.Item(29, 4).Value = $false
$ws.Cells.Item(29, 5).Value = $false
$ws.Cells.Item(29, 6).Value = $true
$ws.Cells.Item(29, 7).Value = $false
$ws.Cells.Item(29, 8).Value = $true
$ws.Cells.Item(29, 9).Value = $false
$ws.Cells.Item(29, 10).Value = 2
$ws.Cells.Item(30, 1).Value = "Discussions viewed"
$ws.Cells.Item(30, 2).Value = $false
$ws.Cells.Item(30, 3).Value = $false
$ws.Cells.Item(30, 4).Value = $false
$ws.Cells.Item(30, 5).Value = $false
$ws.Cells.Item(30, 6).Value = $false
$ws.Cells.Item(30, 7).Value = $false
$ws.Cells.Item(30, 8).Value = $true
$ws.Cells.Item(30, 9).Value = $false
$ws.Cells.Item(30, 10).Value = 1
$ws.Cells.Item(31, 1).Value = "Quizzes started"
$ws.Cells.Item(31, 2).Value = $false
$ws.Cells.Item(31, 3).Value = $false
$ws.Cells.Item(31, 4).Value = $false
$ws.Cells.Item(31, 5).Value = $false
$ws.Cells.Item(31, 6).Value = $false
$ws.Cells.Item(31, 7).Value = $false
$ws.Cells.Item(31, 8).Value = $true
$ws.Cells.Item(31, 9).Value = $false
$ws.Cells.Item(31, 10).Value = 1
$ws.Cells.Item(32, 1).Value = "Clicks on forum"
$ws.Cells.Item(32, 2).Value = $false
$ws.Cells.Item(32, 3).Value = $false
$ws.Cells.Item(32, 4).Value = $false
$ws.Cells.Item(32, 5).Value = $false
$ws.Cells.Item(32, 6).Value = $false
$ws.Cells.Item(32, 7).Value = $false
$ws.Cells.Item(32, 8).Value = $true
$ws.Cells.Item(32, 9).Value = $false
$ws.Cells.Item(32, 10).Value = 1
$ws.Cells.Item(33, 1).Value = "Start of Session 10 (%)"
$ws.Cells.Item(33, 2).Value = $false
$ws.Cells.Item(33, 3).Value = $false
$ws.Cells.Item(33, 4).Value = $false
$ws.Cells.Item(33, 5).Value = $false
$ws.Cells.Item(33, 6).Value = $false
$ws.Cells.Item(33, 7).Value = $false
$ws.Cells.Item(33, 8).Value = $true
$ws.Cells.Item(33, 9).Value = $false
$ws.Cells.Item(33, 10).Value = 1
$ws.Cells.Item(34, 1).Value = "Start of Session 9 (%)"
$ws.Cells.Item(34, 2).Value = $false
$ws.Cells.Item(34, 3).Value = $false
$ws.Cells.Item(34, 4).Value = $false
$ws.Cells.Item(34, 5).Value = $false
$ws.Cells.Item(34, 6).Value = $false
$ws.Cells.Item(34, 7).Value = $false
$ws.Cells.Item(34, 8).Value = $true
$ws.Cells.Item(34, 9).Value = $false
$ws.Cells.Item(34, 10).Value = 1
$ws.Cells.Item(35, 1).Value = "Start of Session 8 (%)"
$ws.Cells.Item(35, 2).Value = $false
$ws.Cells.Item(35, 3).Value = $false
$ws.Cells.Item(35, 4).Value = $false
$ws.Cells.Item(35, 5).Value = $false
$ws.Cells.Item(35, 6).Value = $false
$ws.Cells.Item(35, 7).Value = $false
$ws.Cells.Item(35, 8).Value = $true
$ws.Cells.Item(35, 9).Value = $false
$ws.Cells.Item(35, 10).Value = 1
$ws.Cells.Item(36, 1).Value = "Start of Session 6 (%)"
$ws.Cells.Item(36, 2).Value = $false
$ws.Cells.Item(36, 3).Value = $false
$ws.Cells.Item(36, 4).Value = $false
$ws.Cells.Item(36, 5).Value = $false
$ws.Cells.Item(36, 6).Value = $false
$ws.Cells.Item(36, 7).Value = $false
$ws.Cells.Item(36, 8).Value = $true
$ws.Cells.Item(36, 9).Value = $false
$ws.Cells.Item(36, 10).Value = 1
$ws.Cells.Item(37, 1).Value = "Links viewed"
$ws.Cells.Item(37, 2).Value = $false
$ws.Cells.Item(37, 3).Value = $false
$ws.Cells.Item(37, 4).Value = $false
$ws.Cells.Item(37, 5).Value = $false
$ws.Cells.Item(37, 6).Value = $false
$ws.Cells.Item(37, 7).Value = $false
$ws.Cells.Item(37, 8).Value = $true
$ws.Cells.Item(37, 9).Value = $false
$ws.Cells.Item(37, 10).Value = 1

$ws = $wb.Worksheets.Item("final_gifted")
$ws.Cells.Item(2, 1).Value = "Total time online (min)"
$ws.Cells.Item(2, 2).Value = $true
$ws.Cells.Item(2, 3).Value = $true
$ws.Cells.Item(2, 4).Value = $true
$ws.Cells.Item(2, 5).Value = $true
$ws.Cells.Item(2, 6).Value = $true
$ws.Cells.Item(2, 7).Value = $false
$ws.Cells.Item(2, 8).Value = $true
$ws.Cells.Item(2, 9).Value = $false
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(3, 1).Value = "Average session duration (min)"
$ws.Cells.Item(3, 2).Value = $true
$ws.Cells.Item(3, 3).Value = $true
$ws.Cells.Item(3, 4).Value = $true
$ws.Cells.Item(3, 5).Value = $true
$ws.Cells.Item(3, 6).Value = $true
$ws.Cells.Item(3, 7).Value = $false
$ws.Cells.Item(3, 8).Value = $true
$ws.Cells.Item(3, 9).Value = $false
$ws.Cells.Item(3, 10).Value = 6
$ws.Cells.Item(4, 1).Value = "On/off campus click ratio"
$ws.Cells.Item(4, 2).Value = $true
$ws.Cells.Item(4, 3).Value = $true
$ws.Cells.Item(4, 4).Value = $true
$ws.Cells.Item(4, 5).Value = $true
$ws.Cells.Item(4, 6).Value = $true
$ws.Cells.Item(4, 7).Value = $false
$ws.Cells.Item(4, 8).Value = $true
$ws.Cells.Item(4, 9).Value = $false
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(5, 1).Value = "Clicks (% of course total)"
$ws.Cells.Item(5, 2).Value = $true
$ws.Cells.Item(5, 3).Value = $true
$ws.Cells.Item(5, 4).Value = $true
$ws.Cells.Item(5, 5).Value = $true
$ws.Cells.Item(5, 6).Value = $true
$ws.Cells.Item(5, 7).Value = $false
$ws.Cells.Item(5, 8).Value = $true
$ws.Cells.Item(5, 9).Value = $false
$ws.Cells.Item(5, 10).Value = 6
$ws.Cells.Item(6, 1).Value = "Days with no interaction"
$ws.Cells.Item(6, 2).Value = $true
$ws.Cells.Item(6, 3).Value = $true
$ws.Cells.Item(6, 4).Value = $false
$ws.Cells.Item(6, 5).Value = $true
$ws.Cells.Item(6, 6).Value = $true
$ws.Cells.Item(6, 7).Value = $false
$ws.Cells.Item(6, 8).Value = $true
$ws.Cells.Item(6, 9).Value = $false
$ws.Cells.Item(6, 10).Value = 5
$ws.Cells.Item(7, 1).Value = "Clicks per session"
$ws.Cells.Item(7, 2).Value = $true
$ws.Cells.Item(7, 3).Value = $true
$ws.Cells.Item(7, 4).Value = $false
$ws.Cells.Item(7, 5).Value = $true
$ws.Cells.Item(7, 6).Value = $true
$ws.Cells.Item(7, 7).Value = $false
$ws.Cells.Item(7, 8).Value = $true
$ws.Cells.Item(7, 9).Value = $false
$ws.Cells.Item(7, 10).Value = 5
$ws.Cells.Item(8, 1).Value = "Average grade of assignments"
$ws.Cells.Item(8, 2).Value = $true
$ws.Cells.Item(8, 3).Value = $true
$ws.Cells.Item(8, 4).Value = $true
$ws.Cells.Item(8, 5).Value = $true
$ws.Cells.Item(8, 6).Value = $false
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = $true
$ws.Cells.Item(8, 9).Value = $false
$ws.Cells.Item(8, 10).Value = 5
$ws.Cells.Item(9, 1).Value = "Largest period of inactivity (h)"
$ws.Cells.Item(9, 2).Value = $true
$ws.Cells.Item(9, 3).Value = $true
$ws.Cells.Item(9, 4).Value = $false
$ws.Cells.Item(9, 5).Value = $true
$ws.Cells.Item(9, 6).Value = $true
$ws.Cells.Item(9, 7).Value = $false
$ws.Cells.Item(9, 8).Value = $true
$ws.Cells.Item(9, 9).Value = $false
$ws.Cells.Item(9, 10).Value = 5
$ws.Cells.Item(10, 1).Value = "Start of Session 1 (%)"
$ws.Cells.Item(10, 2).Value = $true
$ws.Cells.Item(10, 3).Value = $true
$ws.Cells.Item(10, 4).Value = $false
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(10, 6).Value = $true
$ws.Cells.Item(10, 7).Value = $false
$ws.Cells.Item(10, 8).Value = $true
$ws.Cells.Item(10, 9).Value = $false
$ws.Cells.Item(10, 10).Value = 5
$ws.Cells.Item(11, 1).Value = "Clicks per day"
$ws.Cells.Item(11, 2).Value = $true
$ws.Cells.Item(11, 3).Value = $false
$ws.Cells.Item(11, 4).Value = $true
$ws.Cells.Item(11, 5).Value = $true
$ws.Cells.Item(11, 6).Value = $true
$ws.Cells.Item(11, 7).Value = $false
$ws.Cells.Item(11, 8).Value = $true
$ws.Cells.Item(11, 9).Value = $false
$ws.Cells.Item(11, 10).Value = 5
$ws.Cells.Item(12, 1).Value = "Resources viewed"
$ws.Cells.Item(12, 2).Value = $true
$ws.Cells.Item(12, 3).Value = $true
$ws.Cells.Item(12, 4).Value = $false
$ws.Cells.Item(12, 5).Value = $true
$ws.Cells.Item(12, 6).Value = $true
$ws.Cells.Item(12, 7).Value = $false
$ws.Cells.Item(12, 8).Value = $true
$ws.Cells.Item(12, 9).Value = $false
$ws.Cells.Item(12, 10).Value = 5
$ws.Cells.Item(13, 1).Value = "Clicks on course"
$ws.Cells.Item(13, 2).Value = $true
$ws.Cells.Item(13, 3).Value = $true
$ws.Cells.Item(13, 4).Value = $false
$ws.Cells.Item(13, 5).Value = $false
$ws.Cells.Item(13, 6).Value = $true
$ws.Cells.Item(13, 7).Value = $false
$ws.Cells.Item(13, 8).Value = $true
$ws.Cells.Item(13, 9).Value = $false
$ws.Cells.Item(13, 10).Value = 4
$ws.Cells.Item(14, 1).Value = "Days with no interaction (%)"
$ws.Cells.Item(14, 2).Value = $true
$ws.Cells.Item(14, 3).Value = $false
$ws.Cells.Item(14, 4).Value = $false
$ws.Cells.Item(14, 5).Value = $true
$ws.Cells.Item(14, 6).Value = $true
$ws.Cells.Item(14, 7).Value = $false
$ws.Cells.Item(14, 8).Value = $true
$ws.Cells.Item(14, 9).Value = $false
$ws.Cells.Item(14, 10).Value = 4
$ws.Cells.Item(15, 1).Value = "Start of Session 2 (%)"
$ws.Cells.Item(15, 2).Value = $true
$ws.Cells.Item(15, 3).Value = $false
$ws.Cells.Item(15, 4).Value = $false
$ws.Cells.Item(15, 5).Value = $true
$ws.Cells.Item(15, 6).Value = $true
$ws.Cells.Item(15, 7).Value = $false
$ws.Cells.Item(15, 8).Value = $true
$ws.Cells.Item(15, 9).Value = $false
$ws.Cells.Item(15, 10).Value = 4
$ws.Cells.Item(16, 1).Value = "Start of Session 4 (%)"
$ws.Cells.Item(16, 2).Value = $true
$ws.Cells.Item(16, 3).Value = $false
$ws.Cells.Item(16, 4).Value = $false
$ws.Cells.Item(16, 5).Value = $false
$ws.Cells.Item(16, 6).Value = $true
$ws.Cells.Item(16, 7).Value = $false
$ws.Cells.Item(16, 8).Value = $true
$ws.Cells.Item(16, 9).Value = $false
$ws.Cells.Item(16, 10).Value = 3
$ws.Cells.Item(17, 1).Value = "Number of clicks"
$ws.Cells.Item(17, 2).Value = $false
$ws.Cells.Item(17, 3).Value = $false
$ws.Cells.Item(17, 4).Value = $false
$ws.Cells.Item(17, 5).Value = $true
$ws.Cells.Item(17, 6).Value = $true
$ws.Cells.Item(17, 7).Value = $false
$ws.Cells.Item(17, 8).Value = $true
$ws.Cells.Item(17, 9).Value = $false
$ws.Cells.Item(17, 10).Value = 3
$ws.Cells.Item(18, 1).Value = "Quizzes started"
$ws.Cells.Item(18, 2).Value = $false
$ws.Cells.Item(18, 3).Value = $false
$ws.Cells.Item(18, 4).Value = $true
$ws.Cells.Item(18, 5).Value = $false
$ws.Cells.Item(18, 6).Value = $false
$ws.Cells.Item(18, 7).Value = $false
$ws.Cells.Item(18, 8).Value = $true
$ws.Cells.Item(18, 9).Value = $false
$ws.Cells.Item(18, 10).Value = 2
$ws.Cells.Item(19, 1).Value = "Submissions (% of course total)"
$ws.Cells.Item(19, 2).Value = $true
$ws.Cells.Item(19, 3).Value = $false
$ws.Cells.Item(19, 4).Value = $false
$ws.Cells.Item(19, 5).Value = $false
$ws.Cells.Item(19, 6).Value = $false
$ws.Cells.Item(19, 7).Value = $false
$ws.Cells.Item(19, 8).Value = $true
$ws.Cells.Item(19, 9).Value = $false
$ws.Cells.Item(19, 10).Value = 2
$ws.Cells.Item(20, 1).Value = "Forum posts"
$ws.Cells.Item(20, 2).Value = $false
$ws.Cells.Item(20, 3).Value = $false
$ws.Cells.Item(20, 4).Value = $true
$ws.Cells.Item(20, 5).Value = $false
$ws.Cells.Item(20, 6).Value = $false
$ws.Cells.Item(20, 7).Value = $false
$ws.Cells.Item(20, 8).Value = $true
$ws.Cells.Item(20, 9).Value = $false
$ws.Cells.Item(20, 10).Value = 2
$ws.Cells.Item(21, 1).Value = "Number of days"
$ws.Cells.Item(21, 2).Value = $true
$ws.Cells.Item(21, 3).Value = $false
$ws.Cells.Item(21, 4).Value = $false
$ws.Cells.Item(21, 5).Value = $false
$ws.Cells.Item(21, 6).Value = $false
$ws.Cells.Item(21, 7).Value = $false
$ws.Cells.Item(21, 8).Value = $true
$ws.Cells.Item(21, 9).Value = $false
$ws.Cells.Item(21, 10).Value = 2
$ws.Cells.Item(22, 1).Value = "Assignments viewed"
$ws.Cells.Item(22, 2).Value = $true
$ws.Cells.Item(22, 3).Value = $false
$ws.Cells.Item(22, 4).Value = $false
$ws.Cells.Item(22, 5).Value = $false
$ws.Cells.Item(22, 6).Value = $false
$ws.Cells.Item(22, 7).Value = $false
$ws.Cells.Item(22, 8).Value = $true
$ws.Cells.Item(22, 9).Value = $false
$ws.Cells.Item(22, 10).Value = 2
$ws.Cells.Item(23, 1).Value = "Clicks on folder"
$ws.Cells.Item(23, 2).Value = $false
$ws.Cells.Item(23, 3).Value = $false
$ws.Cells.Item(23, 4).Value = $true
$ws.Cells.Item(23, 5).Value = $false
$ws.Cells.Item(23, 6).Value = $false
$ws.Cells.Item(23, 7).Value = $false
$ws.Cells.Item(23, 8).Value = $true
$ws.Cells.Item(23, 9).Value = $false
$ws.Cells.Item(23, 10).Value = 2
$ws.Cells.Item(24, 1).Value = "Start of Session 7 (%)"
$ws.Cells.Item(24, 2).Value = $false
$ws.Cells.Item(24, 3).Value = $false
$ws.Cells.Item(24, 4).Value = $false
$ws.Cells.Item(24, 5).Value = $false
$ws.Cells.Item(24, 6).Value = $true
$ws.Cells.Item(24, 7).Value = $false
$ws.Cells.Item(24, 8).Value = $true
$ws.Cells.Item(24, 9).Value = $false
$ws.Cells.Item(24, 10).Value = 2
$ws.Cells.Item(25, 1).Value = "Start of Session 6 (%)"
$ws.Cells.Item(25, 2).Value = $false
$ws.Cells.Item(25, 3).Value = $false
$ws.Cells.Item(25, 4).Value = $false
$ws.Cells.Item(25, 5).Value = $false
$ws.Cells.Item(25, 6).Value = $true
$ws.Cells.Item(25, 7).Value = $false
$ws.Cells.Item(25, 8).Value = $true
$ws.Cells.Item(25, 9).Value = $false
$ws.Cells.Item(25, 10).Value = 2
$ws.Cells.Item(26, 1).Value = "Start of Session 5 (%)"
$ws.Cells.Item(26, 2).Value = $false
$ws.Cells.Item(26, 3).Value = $false
$ws.Cells.Item(26, 4).Value = $false
$ws.Cells.Item(26, 5).Value = $false
$ws.Cells.Item(26, 6).Value = $true
$ws.Cells.Item(26, 7).Value = $false
$ws.Cells.Item(26, 8).Value = $true
$ws.Cells.Item(26, 9).Value = $false
$ws.Cells.Item(26, 10).Value = 2
$ws.Cells.Item(27, 1).Value = "Start of Session 3 (%)"
$ws.Cells.Item(27, 2).Value = $false
$ws.Cells.Item(27, 3).Value = $false
$ws.Cells.Item(27, 4).Value = $false
$ws.Cells.Item(27, 5).Value = $false
$ws.Cells.Item(27, 6).Value = $true
$ws.Cells.Item(27, 7).Value = $false
$ws.Cells.Item(27, 8).Value = $true
$ws.Cells.Item(27, 9).Value = $false
$ws.Cells.Item(27, 10).Value = 2
$ws.Cells.Item(28, 1).Value = "Clicks on campus"
$ws.Cells.Item(28, 2).Value = $false
$ws.Cells.Item(28, 3).Value = $false
$ws.Cells.Item(28, 4).Value = $false
$ws.Cells.Item(28, 5).Value = $false
$ws.Cells.Item(28, 6).Value = $true
$ws.Cells.Item(28, 7).Value = $false
$ws.Cells.Item(28, 8).Value = $true
$ws.Cells.Item(28, 9).Value = $false
$ws.Cells.Item(28, 10).Value = 2
$ws.Cells.Item(29, 1).Value = "Links viewed"
$ws.Cells.Item(29, 2).Value = $false
$ws.Cells.Item(29, 3).Value = $false
$ws.Cells.Item(29, 4).Value = $true
$ws.Cells.Item(29, 5).Value = $false
$ws.Cells.Item(29, 6).Value = $false
$ws.Cells.Item(29, 7).Value = $false
$ws.Cells.Item(29, 8).Value = $true
$ws.Cells.Item(29, 9).Value = $false
$ws.Cells.Item(29, 10).Value = 2
$ws.Cells.Item(30, 1).Value = "Assignments submitted"
$ws.Cells.Item(30, 2).Value = $false
$ws.Cells.Item(30, 3).Value = $false
$ws.Cells.Item(30, 4).Value = $false
$ws.Cells.Item(30, 5).Value = $false
$ws.Cells.Item(30, 6).Value = $false
$ws.Cells.Item(30, 7).Value = $false
$ws.Cells.Item(30, 8).Value = $true
$ws.Cells.Item(30, 9).Value = $false
$ws.Cells.Item(30, 10).Value = 1
$ws.Cells.Item(31, 1).Value = "Discussions viewed"
$ws.Cells.Item(31, 2).Value = $false
$ws.Cells.Item(31, 3).Value = $false
$ws.Cells.Item(31, 4).Value = $false
$ws.Cells.Item(31, 5).Value = $false
$ws.Cells.Item(31, 6).Value = $false
$ws.Cells.Item(31, 7).Value = $false
$ws.Cells.Item(31, 8).Value = $true
$ws.Cells.Item(31, 9).Value = $false
$ws.Cells.Item(31, 10).Value = 1
$ws.Cells.Item(32, 1).Value = "Number of sessions"
$ws.Cells.Item(32, 2).Value = $false
$ws.Cells.Item(32, 3).Value = $false
$ws.Cells.Item(32, 4).Value = $false
$ws.Cells.Item(32, 5).Value = $false
$ws.Cells.Item(32, 6).Value = $false
$ws.Cells.Item(32, 7).Value = $false
$ws.Cells.Item(32, 8).Value = $true
$ws.Cells.Item(32, 9).Value = $false
$ws.Cells.Item(32, 10).Value = 1
$ws.Cells.Item(33, 1).Value = "Clicks on forum"
$ws.Cells.Item(33, 2).Value = $false
$ws.Cells.Item(33, 3).Value = $false
$ws.Cells.Item(33, 4).Value = $false
$ws.Cells.Item(33, 5).Value = $false
$ws.Cells.Item(33, 6).Value = $false
$ws.Cells.Item(33, 7).Value = $false
$ws.Cells.Item(33, 8).Value = $true
$ws.Cells.Item(33, 9).Value = $false
$ws.Cells.Item(33, 10).Value = 1
$ws.Cells.Item(34, 1).Value = "Files downloaded"
$ws.Cells.Item(34, 2).Value = $false
$ws.Cells.Item(34, 3).Value = $false
$ws.Cells.Item(34, 4).Value = $false
$ws.Cells.Item(34, 5).Value = $false
$ws.Cells.Item(34, 6).Value = $false
$ws.Cells.Item(34, 7).Value = $false
$ws.Cells.Item(34, 8).Value = $true
$ws.Cells.Item(34, 9).Value = $false
$ws.Cells.Item(34, 10).Value = 1
$ws.Cells.Item(35, 1).Value = "Start of Session 10 (%)"
$ws.Cells.Item(35, 2).Value = $false
$ws.Cells.Item(35, 3).Value = $false
$ws.Cells.Item(35, 4).Value = $false
$ws.Cells.Item(35, 5).Value = $false
$ws.Cells.Item(35, 6).Value = $false
$ws.Cells.Item(35, 7).Value = $false
$ws.Cells.Item(35, 8).Value = $true
$ws.Cells.Item(35, 9).Value = $false
$ws.Cells.Item(35, 10).Value = 1
$ws.Cells.Item(36, 1).Value = "Start of Session 9 (%)"
$ws.Cells.Item(36, 2).Value = $false
$ws.Cells.Item(36, 3).Value = $false
$ws.Cells.Item(36, 4).Value = $false
$ws.Cells.Item(36, 5).Value = $false
$ws.Cells.Item(36, 6).Value = $false
$ws.Cells.Item(36, 7).Value = $false
$ws.Cells.Item(36, 8).Value = $true
$ws.Cells.Item(36, 9).Value = $false
$ws.Cells.Item(36, 10).Value = 1
$ws.Cells.Item(37, 1).Value = "Start of Session 8 (%)"
$ws.Cells.Item(37, 2).Value = $false
$ws.Cells.Item(37, 3).Value = $false
$ws.Cells.Item(37, 4).Value = $false
$ws.Cells.Item(37, 5).Value = $false
$ws.Cells.Item(37, 6).Value = $false
$ws.Cells.Item(37, 7).Value = $false
$ws.Cells.Item(37, 8).Value = $true
$ws.Cells.Item(37, 9).Value = $false
$ws.Cells.Item(37, 10).Value = 1
